# Clear the weekly maintenance-duration values (columns C:E = "w1","w2","w3")
# on the "Maintenance" sheet back down to 0 - this is the data that gets
# (re)written by the pickle-output step of the generating script.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Maintenance")

# Rows 2-30, 37-41 and 44-79 still carried the old placeholder values
# (150 / 50 / 200) in C/D/E - zero them out.
$ws.Range("C2:E30").Value = 0
$ws.Range("C37:E41").Value = 0
$ws.Range("C44:E79").Value = 0

# Rows 42-43 only had a stale value left in column C (D/E were already 0).
$ws.Range("C42:C43").Value = 0

# Reflect where the user's cursor ended up after the run (bottom of sheet).
$ws.Range("N78").Select()
